$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing text storage (no numeric/date auto-conversion),
# then reset the cell style back to Normal so no stray style index is introduced.
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "65.168.15"
Set-TextValue "E2" "  -0.64%  "

# Row 3
Set-TextValue "D3" "3.554.26"
Set-TextValue "E3" "  -0.38%  "

# Row 4
Set-TextValue "E4" "  -0.05%  "

# Row 5
Set-TextValue "D5" "598.23"
Set-TextValue "E5" "  -0.15%  "

# Row 6
Set-TextValue "D6" "133.14"
Set-TextValue "E6" "  -5.12%  "

# Row 7
Set-TextValue "D7" "3.554.02"
Set-TextValue "E7" "  -0.38%  "

# Row 9
Set-TextValue "D9" "0.494"
Set-TextValue "E9" "  -0.34%  "

# Row 10
Set-TextValue "E10" "  -2.46%  "

# Row 11
Set-TextValue "D11" "7.11"
Set-TextValue "E11" "  -0.59%  "

# Row 12
Set-TextValue "D12" "0.389"
Set-TextValue "E12" "  -1.11%  "

# Row 13
Set-TextValue "D13" "4.153.33"
Set-TextValue "E13" "  -0.51%  "

# Row 14
Set-TextValue "E14" "  -2.90%  "

# Row 15
Set-TextValue "D15" "26.90"
Set-TextValue "E15" "  -0.66%  "

# Row 16
Set-TextValue "D16" "3.548.45"
Set-TextValue "E16" "  -0.56%  "

# Row 17
Set-TextValue "E17" "  -0.16%  "

# Row 18
Set-TextValue "D18" "65.251.04"
Set-TextValue "E18" "  -0.34%  "

# Row 19
Set-TextValue "D19" "9.90"
Set-TextValue "E19" "  -3.98%  "

# Row 20
Set-TextValue "D20" "14.39"
Set-TextValue "E20" "  +0.87%  "

# Row 21
Set-TextValue "D21" "5.83"
Set-TextValue "E21" "  -0.58%  "

# Row 22
Set-TextValue "D22" "390.79"
Set-TextValue "E22" "  -1.56%  "

# Row 23
Set-TextValue "E23" "  +1.09%  "

# Row 24
Set-TextValue "D24" "3.697.49"
Set-TextValue "E24" "  -0.45%  "

# Row 25
Set-TextValue "D25" "74.15"
Set-TextValue "E25" "  -0.78%  "

# Row 26
Set-TextValue "E26" "  -0.13%  "

# Row 27
Set-TextValue "D27" "0.0000114"
Set-TextValue "E27" "  -2.49%  "

# Row 28
Set-TextValue "D28" "7.79"
Set-TextValue "E28" "  -0.73%  "

# Row 29
Set-TextValue "D29" "1.57"
Set-TextValue "E29" "  +24.66%  "

# Row 30
Set-TextValue "B30" "InternetComputer(DFINITY)"
Set-TextValue "C30" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D30" "8.56"
Set-TextValue "E30" "  +3.33%  "

# Row 31
Set-TextValue "B31" "Binance-PegBSC-USD"
Set-TextValue "C31" "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue "D31" "0.999"
Set-TextValue "E31" "  +0.18%  "

# Row 32
Set-TextValue "D32" "2.28"
Set-TextValue "E32" "  +0.01%  "

# Row 33
Set-TextValue "D33" "3.550.26"
Set-TextValue "E33" "  -1.00%  "

# Row 34
Set-TextValue "D34" "24.04"
Set-TextValue "E34" "  +0.27%  "

# Row 35
Set-TextValue "E35" "  +0.00%  "

# Row 36
Set-TextValue "D36" "0.147"
Set-TextValue "E36" "  -0.28%  "

# Row 37
Set-TextValue "D37" "170.76"
Set-TextValue "E37" "  +1.41%  "

# Row 38
Set-TextValue "D38" "6.93"
Set-TextValue "E38" "  -1.73%  "

# Row 39
Set-TextValue "E39" "  -0.94%  "

# Row 40
Set-TextValue "D40" "5.05"
Set-TextValue "E40" "  +1.02%  "

# Row 41
Set-TextValue "D41" "0.0814"
Set-TextValue "E41" "  +1.27%  "

# Row 42
Set-TextValue "E42" "  -0.96%  "

# Row 43
Set-TextValue "D43" "26.22"
Set-TextValue "E43" "  -1.86%  "

# Row 44
Set-TextValue "E44" "  +4.43%  "

# Row 45
Set-TextValue "D45" "42.99"
Set-TextValue "E45" "  -0.06%  "

# Row 46
Set-TextValue "D46" "0.999"
Set-TextValue "E46" "  -0.10%  "

# Row 47
Set-TextValue "E47" "  -0.12%  "

# Row 48
Set-TextValue "E48" "  -3.06%  "

# Row 49
Set-TextValue "B49" "Cosmos"
Set-TextValue "C49" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D49" "6.91"
Set-TextValue "E49" "  +1.31%  "

# Row 50
Set-TextValue "B50" "Maker"
Set-TextValue "C50" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D50" "2.442.39"
Set-TextValue "E50" "  +0.18%  "
